# Update column G ("K") values on the active sheet for rows 2-30,
# regenerated from Strike# to K values per the commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 3
    3  = 2
    4  = 10
    5  = 4
    6  = 7
    7  = 2
    8  = 4
    9  = 4
    10 = 5
    11 = 2
    12 = 7
    13 = 4
    14 = 3
    15 = 5
    16 = 3
    17 = 5
    18 = 5
    19 = 8
    20 = 5
    21 = 5
    22 = 4
    23 = 6
    24 = 7
    25 = 5
    26 = 5
    27 = 6
    28 = 3
    29 = 5
    30 = 3
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
